# Fruta / hortaliza, semanal
# Inserts this week's new price record for "Zapallo" (Paine, 1a (cosecha))
# right before the current row 424, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 424; rows 424:453 shift down to 425:454.
$ws.Rows.Item(424).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A424").Value2 = 4
$ws.Range("B424").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C424").Value = "Los Lagos"
$ws.Range("D424").Value2 = 45021
$ws.Range("E424").Value2 = 10
$ws.Range("F424").Value2 = 100112045
$ws.Range("G424").Value = "Zapallo"
$ws.Range("H424").Value = "Paine"
$ws.Range("I424").Value = "1a (cosecha)"
$ws.Range("J424").Value2 = 250
$ws.Range("K424").Value2 = 500
$ws.Range("L424").Value2 = 500
$ws.Range("M424").Value2 = 500
$ws.Range("N424").Value = "$/kilo (volumen en unidades)"
$ws.Range("O424").Value = "Región de O'Higgins"
$ws.Range("P424").Value2 = 500
$ws.Range("Q424").Value2 = 1
$ws.Range("R424").Value = "Hortaliza"
